$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rows = $used.Rows.Count

for ($r = 1; $r -le $rows; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    $textA = $cellA.Text

    if ($textA -eq "🟥") {
        $cellA.Value = "📕"
    } elseif ($textA -eq "⬛") {
        $cellA.Value = "📘"
    } elseif ($textA -eq "🟧") {
        $cellA.Value = "📙"
    } elseif ($textA -eq "🟩") {
        $cellA.Value = "📗"
    }

    $cellB = $ws.Cells.Item($r, 2)
    if ($cellB.Text -eq "noir") {
        $cellB.Value = "bleu"
    }
}
